# Update team-specific time-matrix probabilities on Sheet1 with freshly
# computed team-specific values (simulation logic for these values is not
# yet wired up elsewhere in the workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.1764705882352941
$ws.Cells.Item(2, 3).Value = 0.5966386554621849
$ws.Cells.Item(2, 10).Value = 0.02521008403361345
$ws.Cells.Item(2, 16).Value = 0.1176470588235294
$ws.Cells.Item(2, 19).Value = 0.08403361344537816
$ws.Cells.Item(3, 2).Value = 0.01342281879194631
$ws.Cells.Item(3, 3).Value = 0.03355704697986577
$ws.Cells.Item(3, 10).Value = 0.02684563758389262
$ws.Cells.Item(3, 16).Value = 0.7046979865771812
$ws.Cells.Item(3, 19).Value = 0.2214765100671141
$ws.Cells.Item(4, 10).Value = 0.08
$ws.Cells.Item(4, 15).Value = 0.02
$ws.Cells.Item(4, 16).Value = 0.7
$ws.Cells.Item(4, 19).Value = 0.2
$ws.Cells.Item(6, 2).Value = 0.04661016949152542
$ws.Cells.Item(6, 4).Value = 0.008474576271186441
$ws.Cells.Item(6, 6).Value = 0.07627118644067797
$ws.Cells.Item(6, 10).Value = 0.2076271186440678
$ws.Cells.Item(6, 15).Value = 0.01694915254237288
$ws.Cells.Item(6, 17).Value = 0.1398305084745763
$ws.Cells.Item(6, 18).Value = 0.06779661016949153
$ws.Cells.Item(6, 19).Value = 0.4364406779661017
$ws.Cells.Item(7, 2).Value = 0.06432748538011696
$ws.Cells.Item(7, 4).Value = 0.03508771929824561
$ws.Cells.Item(7, 6).Value = 0.04678362573099415
$ws.Cells.Item(7, 10).Value = 0.1578947368421053
$ws.Cells.Item(7, 15).Value = 0.01169590643274854
$ws.Cells.Item(7, 17).Value = 0.1812865497076023
$ws.Cells.Item(7, 18).Value = 0.06432748538011696
$ws.Cells.Item(7, 19).Value = 0.4385964912280702
$ws.Cells.Item(8, 2).Value = 0.09662921348314607
$ws.Cells.Item(8, 4).Value = 0.01573033707865169
$ws.Cells.Item(8, 6).Value = 0.07191011235955057
$ws.Cells.Item(8, 10).Value = 0.101123595505618
$ws.Cells.Item(8, 15).Value = 0.02022471910112359
$ws.Cells.Item(8, 17).Value = 0.1662921348314607
$ws.Cells.Item(8, 18).Value = 0.1213483146067416
$ws.Cells.Item(8, 19).Value = 0.4067415730337079
$ws.Cells.Item(9, 2).Value = 0.06542056074766354
$ws.Cells.Item(9, 4).Value = 0.02803738317757009
$ws.Cells.Item(9, 6).Value = 0.06074766355140187
$ws.Cells.Item(9, 10).Value = 0.08411214953271028
$ws.Cells.Item(9, 15).Value = 0.004672897196261682
$ws.Cells.Item(9, 17).Value = 0.1775700934579439
$ws.Cells.Item(9, 18).Value = 0.1121495327102804
$ws.Cells.Item(9, 19).Value = 0.4672897196261682
$ws.Cells.Item(10, 2).Value = 0.09682947729220223
$ws.Cells.Item(10, 4).Value = 0.02656383890317052
$ws.Cells.Item(10, 5).Value = 0.000856898029134533
$ws.Cells.Item(10, 6).Value = 0.0805484147386461
$ws.Cells.Item(10, 10).Value = 0.1122536418166238
$ws.Cells.Item(10, 15).Value = 0.01456726649528706
$ws.Cells.Item(10, 17).Value = 0.194515852613539
$ws.Cells.Item(10, 18).Value = 0.09511568123393316
$ws.Cells.Item(10, 19).Value = 0.3787489288774636
$ws.Cells.Item(11, 7).Value = 0.1585365853658537
$ws.Cells.Item(11, 10).Value = 0.06097560975609756
$ws.Cells.Item(11, 11).Value = 0.2032520325203252
$ws.Cells.Item(11, 12).Value = 0.5772357723577236
$ws.Cells.Item(12, 7).Value = 0.7337662337662337
$ws.Cells.Item(12, 10).Value = 0.1883116883116883
$ws.Cells.Item(12, 11).Value = 0.006493506493506494
$ws.Cells.Item(12, 12).Value = 0.05844155844155844
$ws.Cells.Item(12, 19).Value = 0.01298701298701299
$ws.Cells.Item(13, 7).Value = 0.5897435897435898
$ws.Cells.Item(13, 10).Value = 0.3846153846153846
$ws.Cells.Item(13, 19).Value = 0.02564102564102564
$ws.Cells.Item(15, 6).Value = 0.01587301587301587
$ws.Cells.Item(15, 8).Value = 0.1375661375661376
$ws.Cells.Item(15, 10).Value = 0.3650793650793651
$ws.Cells.Item(15, 11).Value = 0.07936507936507936
$ws.Cells.Item(15, 15).Value = 0.06878306878306878
$ws.Cells.Item(15, 19).Value = 0.2222222222222222
$ws.Cells.Item(16, 6).Value = 0.01840490797546012
$ws.Cells.Item(16, 8).Value = 0.1533742331288344
$ws.Cells.Item(16, 9).Value = 0.0736196319018405
$ws.Cells.Item(16, 10).Value = 0.4785276073619632
$ws.Cells.Item(16, 11).Value = 0.1104294478527607
$ws.Cells.Item(16, 13).Value = 0.01840490797546012
$ws.Cells.Item(16, 14).Value = 0.006134969325153374
$ws.Cells.Item(16, 15).Value = 0.049079754601227
$ws.Cells.Item(16, 19).Value = 0.09202453987730061
$ws.Cells.Item(17, 6).Value = 0.03535353535353535
$ws.Cells.Item(17, 8).Value = 0.196969696969697
$ws.Cells.Item(17, 9).Value = 0.08585858585858586
$ws.Cells.Item(17, 10).Value = 0.4267676767676767
$ws.Cells.Item(17, 11).Value = 0.08585858585858586
$ws.Cells.Item(17, 13).Value = 0.01262626262626263
$ws.Cells.Item(17, 15).Value = 0.05303030303030303
$ws.Cells.Item(17, 19).Value = 0.1035353535353535
$ws.Cells.Item(18, 6).Value = 0.02325581395348837
$ws.Cells.Item(18, 8).Value = 0.1627906976744186
$ws.Cells.Item(18, 9).Value = 0.1023255813953488
$ws.Cells.Item(18, 10).Value = 0.4558139534883721
$ws.Cells.Item(18, 11).Value = 0.05116279069767442
$ws.Cells.Item(18, 13).Value = 0.02325581395348837
$ws.Cells.Item(18, 15).Value = 0.05116279069767442
$ws.Cells.Item(18, 19).Value = 0.1302325581395349
$ws.Cells.Item(19, 6).Value = 0.01554828150572831
$ws.Cells.Item(19, 8).Value = 0.2356792144026187
$ws.Cells.Item(19, 9).Value = 0.102291325695581
$ws.Cells.Item(19, 10).Value = 0.3477905073649755
$ws.Cells.Item(19, 11).Value = 0.09656301145662848
$ws.Cells.Item(19, 13).Value = 0.0220949263502455
$ws.Cells.Item(19, 15).Value = 0.0630114566284779
$ws.Cells.Item(19, 19).Value = 0.1170212765957447
